$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old range entirely so stale cells (e.g. old A1) are removed
$ws.Range("A1:D5").Clear()

# New header row (row 1): A1 empty, B1/C1/D1 hold the moved headers
$ws.Range("B1").Value = "腐蚀率(mm/a)"
$ws.Range("C1").Value = "埋存时间(年)"
$ws.Range("D1").Value = "站点"

# Data rows 2-5
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 3.719
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = "沈阳站"

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 3.467
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = "沈阳站"

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 3.215
$ws.Range("C4").Value = 3
$ws.Range("D4").Value = "沈阳站"

$ws.Range("A5").Value = 3
$ws.Range("B5").Value = 2.987
$ws.Range("C5").Value = 4
$ws.Range("D5").Value = "沈阳站"

# Column widths: B -> 13, C -> 13.5546875 (stored OOXML width includes Excel's
# built-in ~5px padding and is rounded to whole pixels, so the ColumnWidth
# value has to be chosen so it lands on the desired stored width after that
# conversion)
$ws.Columns.Item(2).ColumnWidth = 12.25
$ws.Columns.Item(3).ColumnWidth = 12.86

# Selection moves to A6
$ws.Range("A6").Select()
